$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete trailing rows 51-56 (content consolidated/removed)
$ws.Range("A51:A56").EntireRow.Delete()

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "956"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "956"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "1525"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "257"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1782"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "414"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "1360"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "268"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "885"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "929"
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = "78"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "1915000(Nineteen`rLakh Fifteen Thousand)"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "896"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "827"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "743"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "461"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "1500000(Fifteen Lakh)"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "280"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "910"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "870"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "736"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "356"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "1769000(Seventeen`rLakh Sixty Nine`rThousand)"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "294"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "956"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "953"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "772"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "583"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "1529000(Fifteen Lakh`rTwenty Nine Thousand)"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2136"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "34"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "190"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "208"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232"
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "Library ( Books, Journals and e-Resources only)"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "165902000 (Sixteen Crore Fifty Nine Lakh Two Thousand)"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "169511000 (Sixteen Crore Ninety Five Lakh Eleven Thousand)"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153018000 (Fifteen Crore Thirty Lakh Eighteen Thousand)"
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "New Equipment and software for Laboratories"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "634400900 (Sixty Three Crore Forty Four Lakh Nine Hundred)"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "870407000 (Eighty Seven Crore Four Lakh Seven Thousand)"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "685769000 (Sixty Eight Crore Fifty Seven Lakh Sixty Nine Thousand)"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "107318500 (Ten Crore Seventy Three Lakh Eighteen`rThousand Five Hundred )"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "102500000 (Ten Crore Twenty Five Lakh)"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "98500000 (Nine Crore Eighty Five Lakh)"
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "Other expenditure on creation of Capital Assets (For setting up`rclassrooms, seminar hall, conference hall , library, Lab, Engg`rworkshops excluding expenditure on Land and Building)"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "137118500 (Thirteen Crore Seventy One Lakh Eighteen`rThousand Five Hundred )"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "103500000 (Ten Crore Thirty Five Lakh)"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "73400000 (Seven Crore Thirty Four Lakh)"
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = "Financial Year"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "2021-22"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "2020-21"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2019-20"
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = ""
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Utilised Amount"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "Utilised Amount"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "Utilised Amount"
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "Annual Operational Expenditure"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = ""
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = ""
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = ""
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "Salaries (Teaching and Non Teaching staff)"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "3132324000 (Three Hundred Thirteen Crore Twenty Three`rLakh Twenty Four Thousand )"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "3072400000 (Three Hundred Seven Crore Twenty Four Lakh)"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2682258000 (Two Hundred Sixty Eight Crore Twenty Two Lakh`rFifty Eight Thousand)"
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "Maintenance of Academic Infrastructure or consumables and`rother running expenditures(excluding maintenance of hostels`rand allied services,rent of the building, depreciation cost, etc)"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "3931518000 (Three Hundred Ninety Three Crore Fifteen Lakh`rEighteen Thousand)"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "3972820000 (Three Hundred Ninety Seven Crore Twenty Eight`rLakh Twenty Thousand)"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4205644000 (Four Hundred Twenty Crore Fifty Six Lakh Forty`rFour Thousand)"
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "Seminars/Conferences/Workshops"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "46800000 (Four Crore Sixty Eight Lakh)"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "52100000 (Five Crore Twenty One Lakh)"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51900000 (Five Crore Nineteen Lakh)"
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "Financial Year"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "2021-22"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "2020-21"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2019-20"
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "Total no. of Sponsored Projects"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "981"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "879"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "862"
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "Total no. of Funding Agencies"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "218"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "150"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "143"
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "Total Amount Received (Amount in Rupees)"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "2255285600"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "2612279000"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2389745000"
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "Amount Received in Words"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Two Hundred Twenty Five Crore Fifty Two Lakh Eighty Five`rThousand Six Hundred"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "Two Hundred And Sixty One Crore Twenty Two Lakh Seventy`rNine Thousand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "Two Hundred Thirty Eight Crore Ninety Seven Lakh Forty Five`rThousand"
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "Financial Year"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "2021-22"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "2020-21"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2019-20"
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "Total no. of Consultancy Projects"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "381"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "375"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "315"
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "Total no. of Client Organizations"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "198"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "126"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "191"
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "Total Amount Received (Amount in Rupees)"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "356079800"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "241883900"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "310320800"
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "Amount Received in Words"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Thirty Five Crore Sixty Lakh Seventy Nine Thousand Eight`rHundred"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "Twenty Four Crore Eighteen Lakh Eighty Three Thousand Nine`rHundred"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "Thirty One Crore Three Lakh Twenty Thousand Eight Hundred"
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "1. Do your institution buildings have Lifts/Ramps?"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Yes, more than 80% of the buildings"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = ""
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = ""
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "2. Do your institution have provision for walking aids, including wheelchairs and transportation from one building to another for`rhandicapped students?"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Yes"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = ""
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = ""
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "3. Do your institution buildings have specially designed toilets for handicapped students?"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Yes, more than 80% of the buildings"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = ""
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = ""
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "Number of faculty members entered"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "572"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = ""
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = ""
